$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'AU-4,AU-4 (1)'
$ws.Range('A3').Value = 'CM-6 b,SC-5,SC-5 (2)'
$ws.Range('A4').Value = 'AU-7 a,AC-6 (9),AU-7 b,AU-12 (3),AU-8 b,AC-6 (8),CM-5 (1)'
$ws.Range('A5').Value = 'CM-7 b,CM-6 b,AC-17 (9),AC-17 (1)'
$ws.Range('A10').Value = 'CM-7 (5) (b),CM-7 (2)'
$ws.Range('A15').Value = 'IA-8,IA-2,AU-3 (1)'
$ws.Range('A17').Value = 'AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A19').Value = 'CM-6 b,IA-5 (1) (a),IA-5 (1) (b)'
$ws.Range('A21').Value = 'AC-12,SC-10,MA-4 e,MA-4 (7)'
$ws.Range('A22').Value = 'AU-7 a,AU-12 a,CM-6 b,AU-7 (1),AU-3 (1),AU-6 (4),MA-4 (1) (a),CM-5 (1),AU-14 (1),AU-3'
$ws.Range('A25').Value = 'AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A31').Value = 'AU-12 c,AU-12 a,AC-2 (4),AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A34').Value = 'AC-11 b,AC-11 a'
$ws.Range('A42').Value = 'SC-28,SC-28 (1)'
$ws.Range('A45').Value = 'AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 b'
$ws.Range('A50').Value = 'CM-6 b,IA-2 (5)'
$ws.Range('A53').Value = 'SC-13,MA-4 (6)'
$ws.Range('A56').Value = 'AU-12 c,MA-4 (1) (a)'
$ws.Range('A63').Value = 'AU-5 a,AU-5 (1)'
$ws.Range('A67').Value = 'AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A69').Value = 'AU-7 a,AU-12 c,AU-12 a,CM-6 b,AU-7 b,AU-12 (3),AU-8 b,CM-5 (1)'
$ws.Range('A77').Value = 'AU-12 c,AU-12 a,AC-2 (4),AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A80').Value = 'IA-2 (2),IA-2 (1),IA-2 (4),IA-2 (3)'
$ws.Range('A81').Value = 'CM-6 b,CM-5 (3)'
$ws.Range('A86').Value = 'AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A88').Value = 'AU-12 c,AC-6 (9),AC-2 (4),CM-5 (1)'
$ws.Range('A89').Value = 'IA-2 (2),IA-2,IA-2 (5),IA-2 (4),IA-2 (3)'
$ws.Range('A97').Value = 'AU-8 (1) (a),AU-8 (1) (b),AU-8 b'
$ws.Range('A102').Value = 'AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A111').Value = 'AU-5 b,AU-5 a'
$ws.Range('A119').Value = 'AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A124').Value = 'AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A125').Value = 'AC-18 (1),CM-7 a'
$ws.Range('A128').Value = 'CM-6 b,IA-5 (1) (c),CM-7 a'
$ws.Range('A148').Value = 'AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-14 (1),AU-3'
$ws.Range('A157').Value = 'AU-12 c,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-3'
$ws.Range('A181').Value = 'CM-6 b,SC-3'
